$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 512; existing rows 512-638 shift down to 513-639,
# pushing the sheet's dimension from A1:R638 to A1:R639.
$ws.Rows.Item(512).Insert()

# Populate the newly inserted row 512 with the new data record.
$ws.Range("A512").Value = 5
$ws.Range("B512").Value = "Macroferia Regional de Talca"
$ws.Range("C512").Value = "Maule"
$ws.Range("D512").Value = 45204
$ws.Range("E512").Value = 7
$ws.Range("F512").Value = 100112023
$ws.Range("G512").Value = "Brócoli"
$ws.Range("H512").Value = "Sin especificar"
$ws.Range("I512").Value = "Primera"
$ws.Range("J512").Value = 4000
$ws.Range("K512").Value = 1000
$ws.Range("L512").Value = 1000
$ws.Range("M512").Value = 1000
$ws.Range("N512").Value = "`$/unidad"
$ws.Range("O512").Value = "Región del Maule"
$ws.Range("P512").Value = 1000
$ws.Range("Q512").Value = 1
$ws.Range("R512").Value = "Hortaliza"
